# Manage news test cases
$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("LoginPage")
$homeSheet  = $wb.Worksheets.Item("HomePage")

# Keep HomePage's own selection parked at C2 (its tab is no longer the active one).
$homeSheet.Activate()
$homeSheet.Range("C2").Select()

# Add a new regression test row to LoginPage: admin1 / admin
$loginSheet.Range("A7").Value = "admin1"
$loginSheet.Range("B7").Value = "admin"

# LoginPage becomes the active tab again, with the selection parked past the new row.
$loginSheet.Activate()
$loginSheet.Range("E7").Select()
